$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (column I) values using NOAA temperature data
$ws.Range("I2").Value = -1.226851851851833
$ws.Range("I3").Value = 18.89814814814816
$ws.Range("I4").Value = 18.89814814814816

# Updated worst_ashp_cop (column N) values recalculated from updated temperatures
$ws.Range("N2").Value = 1.203236793039155
$ws.Range("N3").Value = 1.349021684597804

# Updated best_ashp_cop (column O) values recalculated from updated temperatures
$ws.Range("O2").Value = 1.257328254301852
$ws.Range("O3").Value = 1.419425825968325
